$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the numeric-column formatting (style) from the last existing data
# column (AM) onto the new column (AN) before writing values into it.
$ws.Range("AM2:AM11").Copy() | Out-Null
$ws.Range("AN2:AN11").PasteSpecial(-4122) | Out-Null

# New column AN: header + values for rows 2-11
$ws.Range("AN1").Value = "02-ago"
$ws.Range("AN2").Value = 15
$ws.Range("AN3").Value = 16
$ws.Range("AN4").Value = 12
$ws.Range("AN5").Value = 12
$ws.Range("AN6").Value = 10
$ws.Range("AN7").Value = 17
$ws.Range("AN8").Value = 13
$ws.Range("AN9").Value = 16
$ws.Range("AN10").Value = 15
$ws.Range("AN11").Value = 14

# Update selection to match the post-edit cursor position
$ws.Range("AN12").Select() | Out-Null
